$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lista cnp")

# Row 5, column C currently holds "Constantin " (with a trailing space).
# Update it to "Constantin" (trailing space removed).
$ws.Range("C5").Value = "Constantin"

# Leave the active selection on C5, matching the saved workbook view state.
$ws.Range("C5").Select()
